$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.038.17'
$ws.Range("E2").Value = '  -7.52%  '
$ws.Range("D3").Value = '1.417.15'
$ws.Range("E3").Value = '  -7.81%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9973'
$ws.Range("E4").Value = '  -0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9987'
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '275.16'
$ws.Range("E6").Value = '  -5.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3693'
$ws.Range("E7").Value = '  -5.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3129'
$ws.Range("E8").Value = '  -2.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.66'
$ws.Range("E9").Value = '  -8.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.045'
$ws.Range("E10").Value = '  -2.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06505'
$ws.Range("E11").Value = '  -9.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9997'
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.483'
$ws.Range("E13").Value = '  -5.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.74'
$ws.Range("E14").Value = '  -3.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.201'
$ws.Range("E15").Value = '  -6.33%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001020'
$ws.Range("E16").Value = '  -6.70%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.410.65'
$ws.Range("E17").Value = '  -8.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.05690'
$ws.Range("E18").Value = '  -14.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9991'
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.21'
$ws.Range("E20").Value = '  -15.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.608'
$ws.Range("E21").Value = '  -8.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.79'
$ws.Range("E22").Value = '  -4.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.10'
$ws.Range("E23").Value = '  +2.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.276'
$ws.Range("E24").Value = '  -4.57%  '
$ws.Range("D25").Value = '20.014.58'
$ws.Range("E25").Value = '  -7.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.282'
$ws.Range("E26").Value = '  -4.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '135.50'
$ws.Range("E27").Value = '  -10.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.10'
$ws.Range("E28").Value = '  -7.50%  '
$ws.Range("D29").Value = '1.570.11'
$ws.Range("E29").Value = '  -8.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '109.83'
$ws.Range("E30").Value = '  -6.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.029'
$ws.Range("E31").Value = '  -17.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.355'
$ws.Range("E32").Value = '  -11.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8232'
$ws.Range("E33").Value = '  -13.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07714'
$ws.Range("E34").Value = '  -4.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.451'
$ws.Range("E35").Value = '  -0.84%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05884'
$ws.Range("E36").Value = '  -0.56%  '
$ws.Range("B37").Value = 'WEMIXTOKEN'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.454'
$ws.Range("E37").Value = '  -2.32%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.896'
$ws.Range("E38").Value = '  -5.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9984'
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02075'
$ws.Range("E40").Value = '  -6.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.59'
$ws.Range("E41").Value = '  -6.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1920'
$ws.Range("E42").Value = '  -5.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.092'
$ws.Range("E43").Value = '  -7.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5307'
$ws.Range("E44").Value = '  -8.56%  '
$ws.Range("E45").Value = '  -6.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.528'
$ws.Range("E46").Value = '  -5.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5170'
$ws.Range("E47").Value = '  -7.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '113.36'
$ws.Range("E48").Value = '  -2.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.767'
$ws.Range("E49").Value = '  -6.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.047'
$ws.Range("E50").Value = '  -9.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9990'
$ws.Range("E51").Value = '  -0.21%  '
